# Shear Center Reference point added
# The shear center reference point must now be added for all sections.
# The data column that used to hold the "Edges" reference (D2) is shifted
# down by one (a new reference-point column is inserted ahead of it),
# so the existing values move from {2,3} to {1,2}.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 2

# Reflect the user's last on-sheet selection after the edit.
$ws.Range("Q9").Select()
